# Apply the latest cryptos data pull: updated Price (column D) and Volume(1h)
# percentage (column E) values for each coin row.
#
# Column D stores the price as plain text (it mixes "37.819.79"-style grouped
# thousands with plain decimals like "234.24"), so any value that still parses as
# a plain number is written with the cell pre-set to Text format - otherwise Excel
# would silently reinterpret it as a number (and could drop a trailing zero).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.819.79"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").Value = "2.085.10"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.24"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.03"
$ws.Range("E7").Value = "  +2.86%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.37%  "

$ws.Range("E10").Value = "  +2.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  +2.79%  "

$ws.Range("D12").Value = "2.393.67"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.28"
$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").Value = "2.103.38"
$ws.Range("E17").Value = "  +1.66%  "

$ws.Range("D18").Value = "37.742.87"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.41"
$ws.Range("E20").Value = "  +2.66%  "

$ws.Range("E22").Value = "  +0.94%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.39"
$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("E27").Value = "  +8.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.02"
$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("E29").Value = "  +0.41%  "

$ws.Range("E30").Value = "  +2.16%  "

$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("E32").Value = "  +3.64%  "

$ws.Range("E33").Value = "  +4.50%  "

$ws.Range("E34").Value = "  +2.11%  "

$ws.Range("E35").Value = "  +1.86%  "

$ws.Range("E36").Value = "  +2.31%  "

$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("E39").Value = "  -3.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0990"
$ws.Range("E40").Value = "  +3.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.64"
$ws.Range("E42").Value = "  +1.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.39"
$ws.Range("E43").Value = "  +5.49%  "

$ws.Range("E44").Value = "  +1.10%  "

$ws.Range("D45").Value = "1.460.05"
$ws.Range("E45").Value = "  -1.66%  "

$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.11"
$ws.Range("E47").Value = "  +6.47%  "

$ws.Range("E48").Value = "  +4.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.42"
$ws.Range("E49").Value = "  +3.08%  "

$ws.Range("E50").Value = "  +2.80%  "

$ws.Range("D51").Value = "2.278.24"
$ws.Range("E51").Value = "  +0.76%  "
